# Fruta / hortaliza, semanal
# Insert 3 new weekly records into the "Mango" sheet immediately before the
# current row 756, shifting the existing rows 756:779 down to 759:782, and
# fill the 3 freshly-inserted rows with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 756, 757, 758 (pushes old 756:779 -> 759:782)
$ws.Range("A756:A758").EntireRow.Insert()

# Common (non-varying) column values for every data row in this block
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100108
$producto    = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria   = "Mango"
$variedad    = "Sin especificar"
$unidad      = "`$/bandeja 4 kilos"
$kilosEq     = 4

function Set-MangoRow($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioMed, $origen, $precioKilo) {
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioMed
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKilo
    $ws.Cells.Item($row, 20).Value = $kilosEq
}

Set-MangoRow 756 44509 "Especial" 706  6000 7000 6323 "Perú" 1581
Set-MangoRow 757 44509 "Primera"  1886 5500 7000 5967 "Perú" 1492
Set-MangoRow 758 44509 "Segunda"  1226 5000 7000 5672 "Perú" 1418
